$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Hunk 1: "... overloading rather than enable_if." gains spell-check
# markers around "enable_if", and loses the _GoBack bookmark that used
# to sit at the end of this paragraph.
# ---------------------------------------------------------------------
$p1 = Find-ParagraphContaining "overloading rather than enable_if"
$rng1 = $p1.Range
$xml1 = '<w:p ' + $wns + ' w:rsidR="008F6042" w:rsidRDefault="008F6042" w:rsidP="008F6042">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Rewrite </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>ReadVector</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>/</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>WriteVector</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>/</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> to use </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>true_type</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>/</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>false_type</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> overloading rather than </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/><w:r><w:t>enable_if</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '</w:p>'
$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Hunk 2: remove the completed "Support floating point parameters and
# return values in remote function caller." list item entirely.
# ---------------------------------------------------------------------
$p2 = Find-ParagraphContaining "Support floating point parameters"
$p2.Range.Delete()

# ---------------------------------------------------------------------
# Hunk 3: "Support non-MSVC compilers in remote function caller" run is
# split into "S" / "upport ..." with the _GoBack bookmark (freed up by
# hunk 1) now sitting between them.
# ---------------------------------------------------------------------
$p3 = Find-ParagraphContaining "non-MSVC"
$rng3 = $p3.Range
$xml3 = '<w:p ' + $wns + ' w:rsidR="00F32F6E" w:rsidRDefault="00F32F6E" w:rsidP="00C063C4">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r w:rsidRPr="00D24CFF"><w:t>S</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r w:rsidRPr="00D24CFF"><w:t>upport non-MSVC compilers in remote function caller</w:t></w:r>' +
  '<w:r w:rsidR="00085EA5" w:rsidRPr="00D24CFF"><w:t xml:space="preserve"> (e.g. in calling convention specification)</w:t></w:r>' +
  '<w:r w:rsidRPr="00D24CFF"><w:t>.</w:t></w:r>' +
  '</w:p>'
$rng3.InsertXML($xml3)
